# Update the existing "Weekly Quantity" / "Monthly Trend" sheets and add a
# new "PO Forecast" sheet with forecasted PO quantities.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Rename the "Requested quantity" header on the existing sheets.
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet as the last (3rd) sheet in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# Header row — copy the bold/bordered header style from the "Weekly
# Quantity" sheet, then overwrite the text with the forecast column names.
$headers = @("ds", "PO_Forecast", "yhat_lower", "yhat_upper")
for ($col = 1; $col -le 4; $col++) {
    $ws1.Range("B1").Copy($ws3.Cells.Item(1, $col))
    $ws3.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Forecast data rows: ds (date), PO_Forecast, yhat_lower, yhat_upper.
$rows = @(
    @(44934.99999999999, 140, 122.990620634742,  155.7435235972706),
    @(44941.99999999999, 120, 103.5045355476569,  136.1600791556133),
    @(44969.99999999999,  40, 24.40286940118628,   57.32431599347342),
    @(44983.99999999999,   0, -16.32522400630367,  16.42464502010566),
    @(44990.99999999999,   0, -36.80380210585945,  -3.614515927723478),
    @(44997.99999999999,   0, -56.69043568481458,  -23.0147713307201),
    @(45004.99999999999,   0, -75.11610553530872,  -42.17111725595974),
    @(45011.99999999999,   0, -95.29637233759667,  -63.33079155025606),
    @(45018.99999999999,   0, -115.1962797944934,  -82.62520122298112),
    @(45025.99999999999,   0, -133.6058839056181,  -103.1510098483156),
    @(45032.99999999999,   0, -154.7668768704098,  -122.3919527376114),
    @(45039.99999999999,   0, -174.3592492672907,  -142.5393354436634)
)

$r = 2
foreach ($row in $rows) {
    # Copy the date-formatted style (from "Order Week" column) into A, then
    # set the actual values for every column in the row.
    $ws1.Range("A2").Copy($ws3.Cells.Item($r, 1))
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
